# "adding total of the page"
# The "ANNEXE RAPPORT - Publics" sheet gains an extra ("total") row that is
# merged together with the row above it, and the workbook is extended with
# two copies of that updated sheet -- one renamed "ANNEXE RAPPORT - CMD"
# (replacing "ANNEXE RAPPORT - Publics") and a brand new
# "ANNEXE RAPPORT - LYC" tab. The first sheet ("ANNEXE RAPPORT - CMR")
# becomes the selected/active tab.

$wb = $excel.ActiveWorkbook

$cmrName = "ANNEXE RAPPORT - CMR"
$publicsName = "ANNEXE RAPPORT - Publics"
$cmdName = "ANNEXE RAPPORT - CMD"
$lycName = "ANNEXE RAPPORT - LYC"

# --- Grow the "Publics" sheet with the new total row (merge A4:D4 -> A4:D5) ---
$publics = $wb.Worksheets.Item($publicsName)
$publics.Range("A4:D4").UnMerge()
$publics.Range("A4:D5").Merge()

# --- Duplicate it to become the new "ANNEXE RAPPORT - LYC" sheet ---
$publics.Copy($null, $publics)
$wb.Worksheets.Item("$publicsName (2)").Name = $lycName
$null = $wb.Worksheets.Item($publicsName).Delete()

# --- Duplicate "LYC" again, placed right after "CMR", to become "CMD" ---
$lyc = $wb.Worksheets.Item($lycName)
$cmr = $wb.Worksheets.Item($cmrName)
$lyc.Copy($null, $cmr)
$wb.Worksheets.Item("$lycName (2)").Name = $cmdName

# --- Restore each sheet's own selection ---
$wb.Worksheets.Item($cmdName).Range("B44").Select()
$wb.Worksheets.Item($lycName).Range("B48").Select()

# --- Make "CMR" (the first tab) the active / selected sheet ---
$wb.Worksheets.Item($cmrName).Activate()
$wb.Worksheets.Item($cmrName).Range("B39").Select()
